$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the first four product rows with the new values
$ws.Range("B4").Value = "TELFAST 180MG 20 F.C. TABS"
$ws.Range("H4").Value = "1:0"
$ws.Range("L4").Value = 80
$ws.Range("N4").Value = "0:2"

$ws.Range("B5").Value = "URSOFALK 250MG 20 CAPS."
$ws.Range("H5").Value = "0:0"
$ws.Range("L5").Value = 122
$ws.Range("N5").Value = "1:0"

$ws.Range("B6").Value = "WELLMETAZONE 0.1% CREAM 40 GM"
$ws.Range("H6").Value = "0:0"
$ws.Range("L6").Value = 56
$ws.Range("N6").Value = "1:0"

$ws.Range("B7").Value = "كريم فاتيكا 125 مل"
$ws.Range("H7").Value = "2:0"
$ws.Range("L7").Value = 50
$ws.Range("N7").Value = "1:0"

# Remove the three rows that are no longer needed (old rows 8, 9 and 10),
# shifting the totals/footer rows up.
$ws.Rows(10).Delete()
$ws.Rows(9).Delete()
$ws.Rows(8).Delete()

# Update the total to reflect the remaining four rows
$ws.Range("K8").Value = 308

# The footer row's height is recalculated slightly by Excel after the edit
$ws.Rows(9).RowHeight = 17.25
